$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("regression")

# Copy formats only from row 13 (header-style row) to row 16, and row 14 (data-style row) to row 17
$ws.Range("A13:H13").Copy()
$ws.Range("A16:H16").PasteSpecial(-4122)
$ws.Range("A14:H14").Copy()
$ws.Range("A17:H17").PasteSpecial(-4122)

# New test case TC06 - SelectHotelPage and room count validation
# Set values in the order the new unique shared strings must appear:
# 40 = 27/02/2026, 41 = TC06, 42 = 28/02/2026
$ws.Range("B16").Value = "Location"
$ws.Range("C16").Value = "Hotels"
$ws.Range("D16").Value = "Room Type"
$ws.Range("E16").Value = "Number of Rooms"
$ws.Range("F16").Value = "Adults per Room"
$ws.Range("G16").Value = "Checkindate"
$ws.Range("H16").Value = "Checkoutdate"

$ws.Range("B17").Value = "Sydney"
$ws.Range("C17").Value = "Hotel Creek"
$ws.Range("D17").Value = "Standard"
$ws.Range("E17").Value = "1 - One"
$ws.Range("F17").Value = "2 - Two"
$ws.Range("G17").Value = "27/02/2026"

$ws.Range("A16").Value = "TC06"
$ws.Range("A17").Value = "TC06"

$ws.Range("H17").Value = "28/02/2026"

$ws.Range("G17").Select() | Out-Null
